$wb = $excel.ActiveWorkbook
$wsTemps = $wb.Worksheets.Item("Temps")
$wsCamions = $wb.Worksheets.Item("Camions")

# --- Temps sheet: add new "Temps_total" column (G) ---
$wsTemps.Range("G1").Value = "Temps_total"
$wsTemps.Range("G2").Value = 0
$wsTemps.Range("G3").Value = 1
$wsTemps.Range("G4").Value = 2
$wsTemps.Range("G5").Value = 3
$wsTemps.Range("G6").Value = 4
$wsTemps.Range("G7").Value = 5
$wsTemps.Range("G8").Value = 6

# --- Camions sheet: update values ---
$wsCamions.Range("C2").Value = 150
$wsCamions.Range("C3").Value = 30
$wsCamions.Range("C4").Value = 20

# --- Selections ---
$wsTemps.Range("J6").Select() | Out-Null
$wsCamions.Range("C4").Select() | Out-Null

# --- Active sheet / tab ---
$wsCamions.Activate() | Out-Null
